# edit.ps1
# Applies the diff: inserts a risk-mitigation table after the
# "Risk Mitigation" section intro paragraph, and rewrites the
# "Risk Monitoring" section body (replacing the placeholder text with
# real content, plus three additional explanatory paragraphs).

$d = $word.ActiveDocument

function Find-ParagraphIndex($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------
# 1. Insert the risk mitigation table right after the paragraph that
#    contains "Explain how each risk will be avoided...".
# ---------------------------------------------------------------

$mitigationPromptText = "Explain how each risk will be avoided. (Keep in mind that some activities may help to mitigate more than one risk while some risks may require multiple mitigation activities.)"

$mitigationIndex = Find-ParagraphIndex($mitigationPromptText)
if ($mitigationIndex -eq -1) {
    throw "Could not locate the Risk Mitigation prompt paragraph"
}
$targetPara = $d.Paragraphs($mitigationIndex)

# Insert one character before the very end of the paragraph (i.e.
# strictly inside its text run) so the engine performs a clean
# split/insert rather than merging with the following empty
# paragraph.
$insertPos = $targetPara.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)

$tableXml = @'
    <w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:tblPr>
        <w:tblW w:w="8620" w:type="dxa"/>
        <w:jc w:val="left"/>
        <w:tblInd w:w="108" w:type="dxa"/>
        <w:tblBorders>
          <w:top w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
          <w:left w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
          <w:bottom w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
          <w:right w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
          <w:insideH w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
          <w:insideV w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
        </w:tblBorders>
        <w:shd w:val="clear" w:color="auto" w:fill="cacaca"/>
        <w:tblLayout w:type="fixed"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="4310"/>
        <w:gridCol w:w="4310"/>
      </w:tblGrid>
      <w:tr>
        <w:tblPrEx>
          <w:shd w:val="clear" w:color="auto" w:fill="000000"/>
        </w:tblPrEx>
        <w:trPr>
          <w:trHeight w:val="445" w:hRule="atLeast"/>
          <w:tblHeader/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="4310"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:left w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:bottom w:val="single" w:color="ffffff" w:sz="24" w:space="0" w:shadow="0" w:frame="0"/>
              <w:right w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="000000"/>
            <w:tcMar>
              <w:top w:type="dxa" w:w="0"/>
              <w:left w:type="dxa" w:w="0"/>
              <w:bottom w:type="dxa" w:w="0"/>
              <w:right w:type="dxa" w:w="0"/>
            </w:tcMar>
            <w:vAlign w:val="top"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:keepNext w:val="0"/>
              <w:keepLines w:val="0"/>
              <w:pageBreakBefore w:val="0"/>
              <w:widowControl w:val="1"/>
              <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
              <w:tabs>
                <w:tab w:val="left" w:pos="1440"/>
                <w:tab w:val="left" w:pos="2880"/>
              </w:tabs>
              <w:suppressAutoHyphens w:val="1"/>
              <w:bidi w:val="0"/>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:left="0" w:right="0" w:firstLine="0"/>
              <w:jc w:val="left"/>
              <w:outlineLvl w:val="0"/>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Cambria" w:cs="Cambria" w:hAnsi="Cambria" w:eastAsia="Cambria"/>
                <w:b w:val="1"/>
                <w:bCs w:val="1"/>
                <w:i w:val="0"/>
                <w:iCs w:val="0"/>
                <w:caps w:val="0"/>
                <w:smallCaps w:val="0"/>
                <w:strike w:val="0"/>
                <w:dstrike w:val="0"/>
                <w:outline w:val="0"/>
                <w:color w:val="ffffff"/>
                <w:spacing w:val="0"/>
                <w:kern w:val="0"/>
                <w:position w:val="0"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:u w:val="none"/>
                <w:vertAlign w:val="baseline"/>
                <w:rtl w:val="0"/>
              </w:rPr>
              <w:t>Risk</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="4310"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:left w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:bottom w:val="single" w:color="ffffff" w:sz="24" w:space="0" w:shadow="0" w:frame="0"/>
              <w:right w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="000000"/>
            <w:tcMar>
              <w:top w:type="dxa" w:w="0"/>
              <w:left w:type="dxa" w:w="0"/>
              <w:bottom w:type="dxa" w:w="0"/>
              <w:right w:type="dxa" w:w="0"/>
            </w:tcMar>
            <w:vAlign w:val="top"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:keepNext w:val="0"/>
              <w:keepLines w:val="0"/>
              <w:pageBreakBefore w:val="0"/>
              <w:widowControl w:val="1"/>
              <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
              <w:tabs>
                <w:tab w:val="left" w:pos="1440"/>
                <w:tab w:val="left" w:pos="2880"/>
              </w:tabs>
              <w:suppressAutoHyphens w:val="1"/>
              <w:bidi w:val="0"/>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:left="0" w:right="0" w:firstLine="0"/>
              <w:jc w:val="left"/>
              <w:outlineLvl w:val="0"/>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Cambria" w:cs="Cambria" w:hAnsi="Cambria" w:eastAsia="Cambria"/>
                <w:b w:val="1"/>
                <w:bCs w:val="1"/>
                <w:i w:val="0"/>
                <w:iCs w:val="0"/>
                <w:caps w:val="0"/>
                <w:smallCaps w:val="0"/>
                <w:strike w:val="0"/>
                <w:dstrike w:val="0"/>
                <w:outline w:val="0"/>
                <w:color w:val="ffffff"/>
                <w:spacing w:val="0"/>
                <w:kern w:val="0"/>
                <w:position w:val="0"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:u w:val="none"/>
                <w:vertAlign w:val="baseline"/>
                <w:rtl w:val="0"/>
              </w:rPr>
              <w:t>How to avoid</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tblPrEx>
          <w:shd w:val="clear" w:color="auto" w:fill="cacaca"/>
        </w:tblPrEx>
        <w:trPr>
          <w:trHeight w:val="912" w:hRule="atLeast"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="4310"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="ffffff" w:sz="24" w:space="0" w:shadow="0" w:frame="0"/>
              <w:left w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:bottom w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:right w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="cacaca"/>
            <w:tcMar>
              <w:top w:type="dxa" w:w="0"/>
              <w:left w:type="dxa" w:w="0"/>
              <w:bottom w:type="dxa" w:w="0"/>
              <w:right w:type="dxa" w:w="0"/>
            </w:tcMar>
            <w:vAlign w:val="top"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:keepNext w:val="0"/>
              <w:keepLines w:val="0"/>
              <w:pageBreakBefore w:val="0"/>
              <w:widowControl w:val="1"/>
              <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
              <w:suppressAutoHyphens w:val="0"/>
              <w:bidi w:val="0"/>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:left="0" w:right="0" w:firstLine="0"/>
              <w:jc w:val="left"/>
              <w:outlineLvl w:val="9"/>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:cs="Cambria" w:hAnsi="Times New Roman" w:eastAsia="Cambria"/>
                <w:b w:val="0"/>
                <w:bCs w:val="0"/>
                <w:i w:val="0"/>
                <w:iCs w:val="0"/>
                <w:caps w:val="0"/>
                <w:smallCaps w:val="0"/>
                <w:strike w:val="0"/>
                <w:dstrike w:val="0"/>
                <w:outline w:val="0"/>
                <w:color w:val="000000"/>
                <w:spacing w:val="0"/>
                <w:kern w:val="0"/>
                <w:position w:val="0"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:u w:val="none" w:color="000000"/>
                <w:vertAlign w:val="baseline"/>
                <w:rtl w:val="0"/>
              </w:rPr>
              <w:t xml:space="preserve">Potential for project to be too complex </w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="4310"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="ffffff" w:sz="24" w:space="0" w:shadow="0" w:frame="0"/>
              <w:left w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:bottom w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:right w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="cacaca"/>
            <w:tcMar>
              <w:top w:type="dxa" w:w="0"/>
              <w:left w:type="dxa" w:w="0"/>
              <w:bottom w:type="dxa" w:w="0"/>
              <w:right w:type="dxa" w:w="0"/>
            </w:tcMar>
            <w:vAlign w:val="top"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:keepNext w:val="0"/>
              <w:keepLines w:val="0"/>
              <w:pageBreakBefore w:val="0"/>
              <w:widowControl w:val="1"/>
              <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
              <w:tabs>
                <w:tab w:val="left" w:pos="1440"/>
                <w:tab w:val="left" w:pos="2880"/>
              </w:tabs>
              <w:suppressAutoHyphens w:val="1"/>
              <w:bidi w:val="0"/>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:left="0" w:right="0" w:firstLine="0"/>
              <w:jc w:val="left"/>
              <w:outlineLvl w:val="0"/>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Cambria" w:cs="Cambria" w:hAnsi="Cambria" w:eastAsia="Cambria"/>
                <w:b w:val="0"/>
                <w:bCs w:val="0"/>
                <w:i w:val="0"/>
                <w:iCs w:val="0"/>
                <w:caps w:val="0"/>
                <w:smallCaps w:val="0"/>
                <w:strike w:val="0"/>
                <w:dstrike w:val="0"/>
                <w:outline w:val="0"/>
                <w:color w:val="000000"/>
                <w:spacing w:val="0"/>
                <w:kern w:val="0"/>
                <w:position w:val="0"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:u w:val="none"/>
                <w:vertAlign w:val="baseline"/>
                <w:rtl w:val="0"/>
              </w:rPr>
              <w:t xml:space="preserve">Stay on top of tasks, make sure everyone keeps up with their goals, if unable to provide resources and help. Also adequate planning at the beginning stage will be help to avoid any later complexity issues. </w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tblPrEx>
          <w:shd w:val="clear" w:color="auto" w:fill="cacaca"/>
        </w:tblPrEx>
        <w:trPr>
          <w:trHeight w:val="672" w:hRule="atLeast"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="4310"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:left w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:bottom w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:right w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="e6e6e6"/>
            <w:tcMar>
              <w:top w:type="dxa" w:w="0"/>
              <w:left w:type="dxa" w:w="0"/>
              <w:bottom w:type="dxa" w:w="0"/>
              <w:right w:type="dxa" w:w="0"/>
            </w:tcMar>
            <w:vAlign w:val="top"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:keepNext w:val="0"/>
              <w:keepLines w:val="0"/>
              <w:pageBreakBefore w:val="0"/>
              <w:widowControl w:val="1"/>
              <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
              <w:suppressAutoHyphens w:val="0"/>
              <w:bidi w:val="0"/>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:left="0" w:right="0" w:firstLine="0"/>
              <w:jc w:val="left"/>
              <w:outlineLvl w:val="9"/>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:cs="Cambria" w:hAnsi="Times New Roman" w:eastAsia="Cambria"/>
                <w:b w:val="0"/>
                <w:bCs w:val="0"/>
                <w:i w:val="0"/>
                <w:iCs w:val="0"/>
                <w:caps w:val="0"/>
                <w:smallCaps w:val="0"/>
                <w:strike w:val="0"/>
                <w:dstrike w:val="0"/>
                <w:outline w:val="0"/>
                <w:color w:val="000000"/>
                <w:spacing w:val="0"/>
                <w:kern w:val="0"/>
                <w:position w:val="0"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:u w:val="none" w:color="000000"/>
                <w:vertAlign w:val="baseline"/>
                <w:rtl w:val="0"/>
              </w:rPr>
              <w:t>Possibility of hardware issues in regard to server</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="4310"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:left w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:bottom w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:right w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="e6e6e6"/>
            <w:tcMar>
              <w:top w:type="dxa" w:w="0"/>
              <w:left w:type="dxa" w:w="0"/>
              <w:bottom w:type="dxa" w:w="0"/>
              <w:right w:type="dxa" w:w="0"/>
            </w:tcMar>
            <w:vAlign w:val="top"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:keepNext w:val="0"/>
              <w:keepLines w:val="0"/>
              <w:pageBreakBefore w:val="0"/>
              <w:widowControl w:val="1"/>
              <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
              <w:tabs>
                <w:tab w:val="left" w:pos="1440"/>
                <w:tab w:val="left" w:pos="2880"/>
              </w:tabs>
              <w:suppressAutoHyphens w:val="1"/>
              <w:bidi w:val="0"/>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:left="0" w:right="0" w:firstLine="0"/>
              <w:jc w:val="left"/>
              <w:outlineLvl w:val="0"/>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Cambria" w:cs="Cambria" w:hAnsi="Cambria" w:eastAsia="Cambria"/>
                <w:b w:val="0"/>
                <w:bCs w:val="0"/>
                <w:i w:val="0"/>
                <w:iCs w:val="0"/>
                <w:caps w:val="0"/>
                <w:smallCaps w:val="0"/>
                <w:strike w:val="0"/>
                <w:dstrike w:val="0"/>
                <w:outline w:val="0"/>
                <w:color w:val="000000"/>
                <w:spacing w:val="0"/>
                <w:kern w:val="0"/>
                <w:position w:val="0"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:u w:val="none"/>
                <w:vertAlign w:val="baseline"/>
                <w:rtl w:val="0"/>
              </w:rPr>
              <w:t xml:space="preserve">Have a member of the team act as a Systems Administrator and make sure that it is maintained well and any issues are dealt with swiftly. </w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tblPrEx>
          <w:shd w:val="clear" w:color="auto" w:fill="cacaca"/>
        </w:tblPrEx>
        <w:trPr>
          <w:trHeight w:val="452" w:hRule="atLeast"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="4310"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:left w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:bottom w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:right w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="cacaca"/>
            <w:tcMar>
              <w:top w:type="dxa" w:w="0"/>
              <w:left w:type="dxa" w:w="0"/>
              <w:bottom w:type="dxa" w:w="0"/>
              <w:right w:type="dxa" w:w="0"/>
            </w:tcMar>
            <w:vAlign w:val="top"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:keepNext w:val="0"/>
              <w:keepLines w:val="0"/>
              <w:pageBreakBefore w:val="0"/>
              <w:widowControl w:val="1"/>
              <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
              <w:suppressAutoHyphens w:val="0"/>
              <w:bidi w:val="0"/>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:left="0" w:right="0" w:firstLine="0"/>
              <w:jc w:val="left"/>
              <w:outlineLvl w:val="9"/>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:cs="Cambria" w:hAnsi="Times New Roman" w:eastAsia="Cambria"/>
                <w:b w:val="0"/>
                <w:bCs w:val="0"/>
                <w:i w:val="0"/>
                <w:iCs w:val="0"/>
                <w:caps w:val="0"/>
                <w:smallCaps w:val="0"/>
                <w:strike w:val="0"/>
                <w:dstrike w:val="0"/>
                <w:outline w:val="0"/>
                <w:color w:val="000000"/>
                <w:spacing w:val="0"/>
                <w:kern w:val="0"/>
                <w:position w:val="0"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:u w:val="none" w:color="000000"/>
                <w:vertAlign w:val="baseline"/>
                <w:rtl w:val="0"/>
              </w:rPr>
              <w:t>Possibility of miscommunication in languages/frameworks</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="4310"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:left w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:bottom w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:right w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="cacaca"/>
            <w:tcMar>
              <w:top w:type="dxa" w:w="0"/>
              <w:left w:type="dxa" w:w="0"/>
              <w:bottom w:type="dxa" w:w="0"/>
              <w:right w:type="dxa" w:w="0"/>
            </w:tcMar>
            <w:vAlign w:val="top"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:keepNext w:val="0"/>
              <w:keepLines w:val="0"/>
              <w:pageBreakBefore w:val="0"/>
              <w:widowControl w:val="1"/>
              <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
              <w:tabs>
                <w:tab w:val="left" w:pos="1440"/>
                <w:tab w:val="left" w:pos="2880"/>
              </w:tabs>
              <w:suppressAutoHyphens w:val="1"/>
              <w:bidi w:val="0"/>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:left="0" w:right="0" w:firstLine="0"/>
              <w:jc w:val="left"/>
              <w:outlineLvl w:val="0"/>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Cambria" w:cs="Cambria" w:hAnsi="Cambria" w:eastAsia="Cambria"/>
                <w:b w:val="0"/>
                <w:bCs w:val="0"/>
                <w:i w:val="0"/>
                <w:iCs w:val="0"/>
                <w:caps w:val="0"/>
                <w:smallCaps w:val="0"/>
                <w:strike w:val="0"/>
                <w:dstrike w:val="0"/>
                <w:outline w:val="0"/>
                <w:color w:val="000000"/>
                <w:spacing w:val="0"/>
                <w:kern w:val="0"/>
                <w:position w:val="0"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:u w:val="none"/>
                <w:vertAlign w:val="baseline"/>
                <w:rtl w:val="0"/>
              </w:rPr>
              <w:t xml:space="preserve">Make sure to meet often and keep each other updated on all progress as often as possible. </w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tblPrEx>
          <w:shd w:val="clear" w:color="auto" w:fill="cacaca"/>
        </w:tblPrEx>
        <w:trPr>
          <w:trHeight w:val="612" w:hRule="atLeast"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="4310"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:left w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:bottom w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:right w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="e6e6e6"/>
            <w:tcMar>
              <w:top w:type="dxa" w:w="0"/>
              <w:left w:type="dxa" w:w="0"/>
              <w:bottom w:type="dxa" w:w="0"/>
              <w:right w:type="dxa" w:w="0"/>
            </w:tcMar>
            <w:vAlign w:val="top"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:keepNext w:val="0"/>
              <w:keepLines w:val="0"/>
              <w:pageBreakBefore w:val="0"/>
              <w:widowControl w:val="1"/>
              <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
              <w:suppressAutoHyphens w:val="0"/>
              <w:bidi w:val="0"/>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:left="0" w:right="0" w:firstLine="0"/>
              <w:jc w:val="left"/>
              <w:outlineLvl w:val="9"/>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:cs="Cambria" w:hAnsi="Times New Roman" w:eastAsia="Cambria"/>
                <w:b w:val="0"/>
                <w:bCs w:val="0"/>
                <w:i w:val="0"/>
                <w:iCs w:val="0"/>
                <w:caps w:val="0"/>
                <w:smallCaps w:val="0"/>
                <w:strike w:val="0"/>
                <w:dstrike w:val="0"/>
                <w:outline w:val="0"/>
                <w:color w:val="000000"/>
                <w:spacing w:val="0"/>
                <w:kern w:val="0"/>
                <w:position w:val="0"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:u w:val="none" w:color="000000"/>
                <w:vertAlign w:val="baseline"/>
                <w:rtl w:val="0"/>
              </w:rPr>
              <w:t xml:space="preserve">Most of the Dev Team is inexperienced with a majority of the languages and frameworks </w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:type="dxa" w:w="4310"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:left w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:bottom w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
              <w:right w:val="single" w:color="ffffff" w:sz="8" w:space="0" w:shadow="0" w:frame="0"/>
            </w:tcBorders>
            <w:shd w:val="clear" w:color="auto" w:fill="e6e6e6"/>
            <w:tcMar>
              <w:top w:type="dxa" w:w="0"/>
              <w:left w:type="dxa" w:w="0"/>
              <w:bottom w:type="dxa" w:w="0"/>
              <w:right w:type="dxa" w:w="0"/>
            </w:tcMar>
            <w:vAlign w:val="top"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:keepNext w:val="0"/>
              <w:keepLines w:val="0"/>
              <w:pageBreakBefore w:val="0"/>
              <w:widowControl w:val="1"/>
              <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
              <w:tabs>
                <w:tab w:val="left" w:pos="1440"/>
                <w:tab w:val="left" w:pos="2880"/>
              </w:tabs>
              <w:suppressAutoHyphens w:val="1"/>
              <w:bidi w:val="0"/>
              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:left="0" w:right="0" w:firstLine="0"/>
              <w:jc w:val="left"/>
              <w:outlineLvl w:val="0"/>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Cambria" w:cs="Cambria" w:hAnsi="Cambria" w:eastAsia="Cambria"/>
                <w:b w:val="0"/>
                <w:bCs w:val="0"/>
                <w:i w:val="0"/>
                <w:iCs w:val="0"/>
                <w:caps w:val="0"/>
                <w:smallCaps w:val="0"/>
                <w:strike w:val="0"/>
                <w:dstrike w:val="0"/>
                <w:outline w:val="0"/>
                <w:color w:val="000000"/>
                <w:spacing w:val="0"/>
                <w:kern w:val="0"/>
                <w:position w:val="0"/>
                <w:sz w:val="18"/>
                <w:szCs w:val="18"/>
                <w:u w:val="none"/>
                <w:vertAlign w:val="baseline"/>
                <w:rtl w:val="0"/>
              </w:rPr>
              <w:t>Make sure to read documentation and stay in constant communication with each other, be willing to provide assistance when necessary.</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
'@

$insertRange.InsertXML($tableXml) | Out-Null

# ---------------------------------------------------------------
# 2. Replace the "Risk Monitoring" placeholder paragraph text with
#    the real explanation, and append three more explanatory
#    paragraphs after it.
# ---------------------------------------------------------------

$monitorPromptText = "Explain how the team will monitor each risk to detect if the risk is becoming a problem. (Keep in mind that some activities may help to monitor more than one risk while some risks may require multiple monitors.)"

$monitorIndex = Find-ParagraphIndex($monitorPromptText)
if ($monitorIndex -eq -1) {
    throw "Could not locate the Risk Monitoring prompt paragraph"
}
$monitorPara = $d.Paragraphs($monitorIndex)

$monitorRange = $monitorPara.Range
$monitorRange.Find.Execute($monitorPromptText, $false, $false, $false, $false, $false, $true, 1, $false, "If we attempt to implement too many features or go over your ability, we will make sure to discuss all plans and changes to the objectives, allowing us a platform to provide each other with more perspective and discussion. ", 2) | Out-Null

# Append the three follow-up paragraphs after it, each using the
# "Body Text" style consistent with the rest of the section.
$d.Paragraphs($monitorIndex).Range.InsertParagraphAfter()
$d.Paragraphs($monitorIndex + 1).Range.Text = "To deal with potential hardware issues, we will have someone who will make sure to stay on top of the current server status and to mitigate any issues as they arise. "

$d.Paragraphs($monitorIndex + 1).Range.InsertParagraphAfter()
$d.Paragraphs($monitorIndex + 2).Range.Text = "To deal with miscommunication, someone will maintain a progress chart that includes detailed descriptions of what is being used for each component. "

$d.Paragraphs($monitorIndex + 2).Range.InsertParagraphAfter()
$d.Paragraphs($monitorIndex + 3).Range.Text = "To deal with inexperience we will make sure everyones code looks good and is well documented. "

Write-Host "Edit applied: tables=" $d.Tables.Count " paragraphs=" $d.Paragraphs.Count
